$d = $word.ActiveDocument

# Remove the trailing "source: https://xkcd.com/882/" caption run that
# followed the xkcd comic image, leaving just the image in its paragraph.
$d.Content.Find.Execute("source: https://xkcd.com/882/", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 2)
